# doc(knitwit-78): adding a cover letter
#
# Slide 7 ("Бизнес-модель"): the paragraph under "Поддержка развития
# платформы" was typed as two separate runs ("...разработки и " +
# "обновления продукта."). Merge them back into a single run with the
# combined text (formatting of the first run is kept).
#
# Slide 8 ("План развития"): the bullet "Добавления подписки на курсы на
# платной основе." is corrected to "Добавление подписки на курсы на
# платной основе." and the leading word is split off into its own run
# (e.g. so it can be styled separately later).

$p = $ppt.ActivePresentation

# --- Slide 7: merge the two runs back into one -----------------------
$s7  = $p.Slides.Item(7)
$shp7 = $s7.Shapes.Item(7)
$tr7 = $shp7.TextFrame.TextRange

$mergedText = "В долгосрочной перспективе поддержка развития платформы будет происходить с помощью добавления платных курсов. Это ускорит процесс разработки и обновления продукта."

$run7 = $tr7.Characters(32, 164)
$run7.Text = $mergedText

# --- Slide 8: split "Добавления " off into its own leading run -------
$s8  = $p.Slides.Item(8)
$shp8 = $s8.Shapes.Item(1)
$tr8 = $shp8.TextFrame.TextRange

$lead8 = $tr8.Characters(154, 11)
$lead8.Text = "Добавление "
